$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 154
$ws.Range("F4").Value = 2079
$ws.Range("F5").Value = 4224
$ws.Range("F6").Value = 550
$ws.Range("D7").Value = "老沪闵路1388号舒也时代广场C栋2层 轮客行轮滑馆(闵行店)"
$ws.Range("F7").Value = 1053
$ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202407/82STI5301720059684987.jpeg"
$ws.Range("F8").Value = 1315
$ws.Range("F10").Value = 376
$ws.Range("F11").Value = 2208
$ws.Range("F12").Value = 399
$ws.Range("F13").Value = 661500
$ws.Range("F14").Value = 1636
$ws.Range("F15").Value = 525
$ws.Range("F17").Value = 672
$ws.Range("F18").Value = 543
$ws.Range("F19").Value = 1278
$ws.Range("F20").Value = 2248
$ws.Range("F21").Value = 1142
$ws.Range("F22").Value = 2700
$ws.Range("F23").Value = 1560
$ws.Range("F24").Value = 824
$ws.Range("F25").Value = 1544
$ws.Range("F27").Value = 1086
$ws.Range("F29").Value = 1087
$ws.Range("F31").Value = 81
$ws.Range("F32").Value = 2024
$ws.Range("F33").Value = 1379
$ws.Range("F34").Value = 576
$ws.Range("F35").Value = 1287
$ws.Range("F36").Value = 2586
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 1146
$ws.Range("F39").Value = 32
$ws.Range("F41").Value = 2590
$ws.Range("F43").Value = 987
$ws.Range("F44").Value = 3142
$ws.Range("F45").Value = 1009
# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 64
$ws.Range("F3").Value = 68
$ws.Range("F9").Value = 109
$ws.Range("F10").Value = 483
$ws.Range("F11").Value = 144662
$ws.Range("F12").Value = 144662
$ws.Range("F22").Value = 133
$ws.Range("F26").Value = 572
$ws.Range("F29").Value = 12
$ws.Range("F31").Value = 353
$ws.Range("F32").Value = 273
$ws.Range("F34").Value = 52
$ws.Range("F35").Value = 52
$ws.Range("F38").Value = 202
# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3132
$ws.Range("F5").Value = 244
$ws.Range("F8").Value = 1190
$ws.Range("F9").Value = 638
$ws.Range("F10").Value = 1601
$ws.Range("F12").Value = 99
$ws.Range("F13").Value = 1911
# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1190
$ws.Range("F4").Value = 638
$ws.Range("F6").Value = 1601
$ws.Range("F8").Value = 154
$ws.Range("F9").Value = 2079
$ws.Range("F10").Value = 99
$ws.Range("F11").Value = 1911
$ws.Range("F12").Value = 4224
$ws.Range("F13").Value = 550
$ws.Range("F14").Value = 1315
$ws.Range("F16").Value = 376
$ws.Range("F17").Value = 2208
$ws.Range("F18").Value = 400
$ws.Range("F19").Value = 661504
$ws.Range("F20").Value = 109
$ws.Range("F21").Value = 483
$ws.Range("F22").Value = 1636
$ws.Range("F23").Value = 144662
$ws.Range("F25").Value = 672
$ws.Range("F26").Value = 543
$ws.Range("F27").Value = 1279
$ws.Range("F28").Value = 2248
$ws.Range("F29").Value = 1142
$ws.Range("F30").Value = 2700
$ws.Range("F31").Value = 1560
$ws.Range("F32").Value = 824
$ws.Range("F34").Value = 1544
$ws.Range("F36").Value = 133
$ws.Range("F37").Value = 1086
$ws.Range("F38").Value = 1087
$ws.Range("F39").Value = 1379
$ws.Range("F40").Value = 1287
$ws.Range("F41").Value = 2586
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 1146
$ws.Range("F44").Value = 353
$ws.Range("F45").Value = 273
$ws.Range("F46").Value = 52
$ws.Range("F47").Value = 2590
$ws.Range("F48").Value = 3142
$ws.Range("F49").Value = 202
$ws.Range("F50").Value = 1009
